$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.882.52'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.631.41'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2568'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06338'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.41'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.264'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').Value = '1.631.46'
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('D14').Value = '1.856.71'
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5496'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.74'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('D17').Value = '0.0₅7631'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('D18').Value = '25.916.52'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.34'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.411'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.850'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.026'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.891'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.14'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1250'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.57%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.745'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.54'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04887'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.231'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.180'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('E34').Value = '  +0.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.375'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8971'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.539'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5496'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.97%  '
$ws.Range('D39').Value = '1.114.65'
$ws.Range('E39').Value = '  -2.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01555'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.001'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.574'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7965'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.49'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.65%  '
$ws.Range('D45').Value = '1.765.09'
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('D46').Value = '0.0₈117'
$ws.Range('E46').Value = '  -6.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4438'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.66'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05131'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.520'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.65%  '
